$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40 (shifts existing rows 40-68 down to 41-69)
$ws.Rows(40).Insert()

# Populate the newly inserted row 40 with the new weekly price record
$ws.Cells.Item(40, 1).Value = 4
$ws.Cells.Item(40, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(40, 3).Value = "Los Lagos"
$ws.Cells.Item(40, 4).Value = 44484
$ws.Cells.Item(40, 5).Value = 10
$ws.Cells.Item(40, 6).Value = 100112052
$ws.Cells.Item(40, 7).Value = "Albahaca"
$ws.Cells.Item(40, 8).Value = "Sin especificar"
$ws.Cells.Item(40, 9).Value = "Primera"
$ws.Cells.Item(40, 10).Value = 120
$ws.Cells.Item(40, 11).Value = 6000
$ws.Cells.Item(40, 12).Value = 6000
$ws.Cells.Item(40, 13).Value = 6000
$ws.Cells.Item(40, 14).Value = "`$/paquete"
$ws.Cells.Item(40, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(40, 16).Value = 6000
$ws.Cells.Item(40, 17).Value = 1
$ws.Cells.Item(40, 18).Value = "Hortaliza"
